$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New QC-reject / warranty-return lines (rows 23-30), matching the format of
# the existing data rows (copy formats for columns A, E and H from row 13;
# B and C are left with default/no explicit style, same as the source rows).

$ws.Range("A13").Copy()
$ws.Range("A23:A30").PasteSpecial(-4122)

$ws.Range("E13").Copy()
$ws.Range("E23:E30").PasteSpecial(-4122)

$ws.Range("H13").Copy()
$ws.Range("H23:H30").PasteSpecial(-4122)

$excel.CutCopyMode = 0

$data = @(
    @(11, "RW-09-8415-2-05-0", "OIL Yellow 3G-05KG/CAN", 1200, "Note for line 11"),
    @(12, "RW-09-8502-2-50-0", "CALCIUM OXIDE-50KG/BAG", 1200, "Note for line 12"),
    @(13, "RW-09-8504-2-25-0", "Kalfain 200M-25KG/BAG", 1200, "Note for line 13"),
    @(14, "RW-09-8507-2-25-0", "CALCIUM CARBONATE CARB 1-HO-25KG/BAG", 1200, "Note for line 14"),
    @(15, "RW-09-8508-2-20-0", "VIGOT 15-20KG/BAG", 250, "Note for line 15"),
    @(16, "RW-09-8509-2-30-0", "WA (BF-200)-30KG/BAG", 250, "Note for line 16"),
    @(17, "RW-09-8701-2-25-0", "SOLBIN A-25KG/BAG", 250, "Note for line 17"),
    @(18, "RW-09-8704-2-10-0", "N2O Pharma-10KG/BAG", 1200, "Note for line 18")
)

# Fill column A (No.) first
$r = 23
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $r = $r + 1
}

# Then columns B (item code) and C (item name) together for every row
$r = 23
foreach ($row in $data) {
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}

# Then column E (quantity)
$r = 23
foreach ($row in $data) {
    $ws.Cells.Item($r, 5).Value = $row[3]
    $r = $r + 1
}

# Finally column H (note) for every row
$r = 23
foreach ($row in $data) {
    $ws.Cells.Item($r, 8).Value = $row[4]
    $r = $r + 1
}

# Match the author's final selection / scroll position: row 23 selected,
# view scrolled down toward the new rows.
$ws.Rows.Item(23).Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 16 | Out-Null
